# Update automatico via Actualizar 04-17-2021 12-26-32
# Shift the "Fecha" (date/time) history down one block and stamp the
# newest block with the current refresh timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D30:D43").Value = 44303.47523460648
$ws.Range("D16:D29").Value = 44303.49667890046
$ws.Range("D2:D15").Value  = 44303.51810427992
